$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 5 (pushes old rows 5-8 down to 6-9,
# carrying all of their values/formatting with them).
$ws.Rows("5:5").Insert()

# Populate the newly inserted row 5 with the new record.
$ws.Range("A5").Value = 10
$ws.Range("B5").Value = "Vega Modelo de Temuco"
$ws.Range("C5").Value = "La Araucanía"
$ws.Range("D5").Value = 44495
$ws.Range("E5").Value = 9
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100104
$ws.Range("H5").Value = "Frutos de pepita"
$ws.Range("I5").Value = 100104005
$ws.Range("J5").Value = "Pera asiática"
$ws.Range("K5").Value = "Hosui"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 24000
$ws.Range("O5").Value = 24000
$ws.Range("P5").Value = 24000
$ws.Range("Q5").Value = "`$/bandeja 10 kilos"
$ws.Range("R5").Value = "China"
$ws.Range("S5").Value = 2400
$ws.Range("T5").Value = 10
